# "Add forget me not" - insert a lone space run at the very start of the
# document's first paragraph, ahead of the existing "Lorem ipsum..." run.
$d = $word.ActiveDocument

$firstPara = $d.Paragraphs(1).Range
$insertionPoint = $d.Range($firstPara.Start, $firstPara.Start)
$insertionPoint.InsertBefore(" ")
